# Auto-generated edit script: updates cached market-price / profit
# snapshot values across the 8 job sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR), matching the upstream scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3517.8518   # H40: 3479.6428 -> 3517.8518
$ws.Cells.Item(40, 9).Value = 3066.6667   # I40: 3004.5 -> 3066.6667
$ws.Cells.Item(40, 10).Value = 3743.4443   # J40: 3743.611 -> 3743.4443
$ws.Cells.Item(40, 11).Value = 3066.6667   # K40: 3004.5 -> 3066.6667
$ws.Cells.Item(40, 12).Value = 3743.4443   # L40: 3743.611 -> 3743.4443
$ws.Cells.Item(40, 13).Value = -2891.6667   # M40: -2829.5 -> -2891.6667
$ws.Cells.Item(40, 14).Value = -4093.4443   # N40: -4093.611 -> -4093.4443
$ws.Cells.Item(69, 8).Value = 38336   # H69: 45169.332 -> 38336
$ws.Cells.Item(69, 9).Value = 32500.5   # I69: 42750.5 -> 32500.5
$ws.Cells.Item(69, 11).Value = 97501.5   # K69: 128251.5 -> 97501.5
$ws.Cells.Item(69, 13).Value = -96627.5   # M69: -127377.5 -> -96627.5
$ws.Cells.Item(70, 8).Value = 2890.0625   # H70: 3265.3076 -> 2890.0625
$ws.Cells.Item(70, 9).Value = 1682   # I70: 1892.7142 -> 1682
$ws.Cells.Item(70, 10).Value = 4443.2856   # J70: 4866.6665 -> 4443.2856
$ws.Cells.Item(70, 11).Value = 5046   # K70: 5678.142599999999 -> 5046
$ws.Cells.Item(70, 12).Value = 13329.8568   # L70: 14599.9995 -> 13329.8568
$ws.Cells.Item(70, 13).Value = -4776   # M70: -5408.142599999999 -> -4776
$ws.Cells.Item(70, 14).Value = -13869.8568   # N70: -15139.9995 -> -13869.8568
$ws.Cells.Item(72, 8).Value = 38336   # H72: 45169.332 -> 38336
$ws.Cells.Item(72, 9).Value = 32500.5   # I72: 42750.5 -> 32500.5
$ws.Cells.Item(72, 11).Value = 292504.5   # K72: 384754.5 -> 292504.5
$ws.Cells.Item(72, 13).Value = -288136.5   # M72: -380386.5 -> -288136.5
$ws.Cells.Item(73, 8).Value = 2890.0625   # H73: 3265.3076 -> 2890.0625
$ws.Cells.Item(73, 9).Value = 1682   # I73: 1892.7142 -> 1682
$ws.Cells.Item(73, 10).Value = 4443.2856   # J73: 4866.6665 -> 4443.2856
$ws.Cells.Item(73, 11).Value = 5046   # K73: 5678.142599999999 -> 5046
$ws.Cells.Item(73, 12).Value = 13329.8568   # L73: 14599.9995 -> 13329.8568
$ws.Cells.Item(73, 13).Value = -4110   # M73: -4742.142599999999 -> -4110
$ws.Cells.Item(73, 14).Value = -15201.8568   # N73: -16471.9995 -> -15201.8568
$ws.Cells.Item(97, 8).Value = 1900   # H97: 5000 -> 1900
$ws.Cells.Item(97, 10).Value = 1900   # J97: 5000 -> 1900
$ws.Cells.Item(97, 12).Value = 5700   # L97: 15000 -> 5700
$ws.Cells.Item(97, 14).Value = -6692   # N97: -15992 -> -6692
$ws.Cells.Item(112, 8).Value = 2105.111   # H112: 2085.6 -> 2105.111
$ws.Cells.Item(112, 10).Value = 2124.75   # J112: 2100.889 -> 2124.75
$ws.Cells.Item(112, 12).Value = 6374.25   # L112: 6302.667 -> 6374.25
$ws.Cells.Item(112, 14).Value = -8590.25   # N112: -8518.667000000001 -> -8590.25
$ws.Cells.Item(116, 8).Value = 5254.95   # H116: 5373.6313 -> 5254.95
$ws.Cells.Item(116, 9).Value = 4686   # I116: 4815.6924 -> 4686
$ws.Cells.Item(116, 11).Value = 4686   # K116: 4815.6924 -> 4686
$ws.Cells.Item(116, 13).Value = -1244   # M116: -1373.6924 -> -1244
$ws.Cells.Item(132, 8).Value = 2840.1738   # H132: 2827.7368 -> 2840.1738
$ws.Cells.Item(132, 9).Value = 2840.1738   # I132: 2827.7368 -> 2840.1738
$ws.Cells.Item(132, 11).Value = 8520.5214   # K132: 8483.2104 -> 8520.5214
$ws.Cells.Item(132, 13).Value = -5990.5214   # M132: -5953.2104 -> -5990.5214
$ws.Cells.Item(137, 8).Value = 14846.761   # H137: 14670.083 -> 14846.761
$ws.Cells.Item(137, 9).Value = 32499.783   # I137: 31234.209 -> 32499.783
$ws.Cells.Item(137, 11).Value = 97499.349   # K137: 93702.62699999999 -> 97499.349
$ws.Cells.Item(137, 13).Value = -94949.349   # M137: -91152.62699999999 -> -94949.349
$ws.Cells.Item(138, 8).Value = 5481.0337   # H138: 5482.793 -> 5481.0337
$ws.Cells.Item(138, 9).Value = 6887.3335   # I138: 7641.5 -> 6887.3335
$ws.Cells.Item(138, 11).Value = 20662.0005   # K138: 22924.5 -> 20662.0005
$ws.Cells.Item(138, 13).Value = -15522.0005   # M138: -17784.5 -> -15522.0005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1526.4286   # H2: 1697.5 -> 1526.4286
$ws.Cells.Item(2, 9).Value = 1526.4286   # I2: 1697.5 -> 1526.4286
$ws.Cells.Item(2, 11).Value = 1526.4286   # K2: 1697.5 -> 1526.4286
$ws.Cells.Item(2, 13).Value = -1413.4286   # M2: -1584.5 -> -1413.4286
$ws.Cells.Item(32, 8).Value = 6764747.5   # H32: 6586728.5 -> 6764747.5
$ws.Cells.Item(32, 9).Value = 7148974   # I32: 6950392.5 -> 7148974
$ws.Cells.Item(32, 11).Value = 7148974   # K32: 6950392.5 -> 7148974
$ws.Cells.Item(32, 13).Value = -7148687   # M32: -6950105.5 -> -7148687
$ws.Cells.Item(45, 8).Value = 2282.6843   # H45: 2263.5652 -> 2282.6843
$ws.Cells.Item(45, 9).Value = 1859.6   # I45: 1998.6666 -> 1859.6
$ws.Cells.Item(45, 10).Value = 2433.7856   # J45: 2357.0588 -> 2433.7856
$ws.Cells.Item(45, 11).Value = 1859.6   # K45: 1998.6666 -> 1859.6
$ws.Cells.Item(45, 12).Value = 2433.7856   # L45: 2357.0588 -> 2433.7856
$ws.Cells.Item(45, 13).Value = -1482.6   # M45: -1621.6666 -> -1482.6
$ws.Cells.Item(45, 14).Value = -3187.7856   # N45: -3111.0588 -> -3187.7856
$ws.Cells.Item(74, 8).Value = 13167665   # H74: 13899105 -> 13167665
$ws.Cells.Item(74, 9).Value = 35717610   # I74: 41670250 -> 35717610
$ws.Cells.Item(74, 11).Value = 35717610   # K74: 41670250 -> 35717610
$ws.Cells.Item(74, 13).Value = -35716736   # M74: -41669376 -> -35716736
$ws.Cells.Item(77, 8).Value = 13167665   # H77: 13899105 -> 13167665
$ws.Cells.Item(77, 9).Value = 35717610   # I77: 41670250 -> 35717610
$ws.Cells.Item(77, 11).Value = 178588050   # K77: 208351250 -> 178588050
$ws.Cells.Item(77, 13).Value = -178583682   # M77: -208346882 -> -178583682
$ws.Cells.Item(113, 8).Value = 110000   # H113: 112000 -> 110000
$ws.Cells.Item(113, 10).Value = 110000   # J113: 112000 -> 110000
$ws.Cells.Item(113, 12).Value = 110000   # L113: 112000 -> 110000
$ws.Cells.Item(113, 14).Value = -118678   # N113: -120678 -> -118678
$ws.Cells.Item(116, 8).Value = 1526.4286   # H116: 1697.5 -> 1526.4286
$ws.Cells.Item(116, 9).Value = 1526.4286   # I116: 1697.5 -> 1526.4286
$ws.Cells.Item(116, 11).Value = 1526.4286   # K116: 1697.5 -> 1526.4286
$ws.Cells.Item(116, 13).Value = 767.5714   # M116: 596.5 -> 767.5714
$ws.Cells.Item(122, 8).Value = 3436.2917   # H122: 3542.2173 -> 3436.2917
$ws.Cells.Item(122, 9).Value = 1470.1   # I122: 1522.3334 -> 1470.1
$ws.Cells.Item(122, 11).Value = 4410.299999999999   # K122: 4567.0002 -> 4410.299999999999
$ws.Cells.Item(122, 13).Value = -1960.299999999999   # M122: -2117.0002 -> -1960.299999999999
$ws.Cells.Item(132, 8).Value = 5129.5264   # H132: 5761.4375 -> 5129.5264
$ws.Cells.Item(132, 9).Value = 1592.5   # I132: 1751.4166 -> 1592.5
$ws.Cells.Item(132, 10).Value = 15033.2   # J132: 17791.5 -> 15033.2
$ws.Cells.Item(132, 11).Value = 4777.5   # K132: 5254.2498 -> 4777.5
$ws.Cells.Item(132, 12).Value = 45099.60000000001   # L132: 53374.5 -> 45099.60000000001
$ws.Cells.Item(132, 13).Value = -2247.5   # M132: -2724.2498 -> -2247.5
$ws.Cells.Item(132, 14).Value = -50159.60000000001   # N132: -58434.5 -> -50159.60000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1526.4286   # H3: 1697.5 -> 1526.4286
$ws.Cells.Item(3, 9).Value = 1526.4286   # I3: 1697.5 -> 1526.4286
$ws.Cells.Item(3, 11).Value = 1526.4286   # K3: 1697.5 -> 1526.4286
$ws.Cells.Item(3, 13).Value = -1412.4286   # M3: -1583.5 -> -1412.4286
$ws.Cells.Item(20, 8).Value = 1955.5769   # H20: 1912.7407 -> 1955.5769
$ws.Cells.Item(20, 9).Value = 1495.7368   # I20: 1460.9 -> 1495.7368
$ws.Cells.Item(20, 11).Value = 1495.7368   # K20: 1460.9 -> 1495.7368
$ws.Cells.Item(20, 13).Value = -1248.7368   # M20: -1213.9 -> -1248.7368

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 561191.5600000001   # H31: 713937.9399999999 -> 561191.5600000001
$ws.Cells.Item(31, 9).Value = 9965.842000000001   # I31: 12350.134 -> 9965.842000000001
$ws.Cells.Item(31, 10).Value = 1016551.94   # J31: 1298594.4 -> 1016551.94
$ws.Cells.Item(31, 11).Value = 9965.842000000001   # K31: 12350.134 -> 9965.842000000001
$ws.Cells.Item(31, 12).Value = 1016551.94   # L31: 1298594.4 -> 1016551.94
$ws.Cells.Item(31, 13).Value = -9670.842000000001   # M31: -12055.134 -> -9670.842000000001
$ws.Cells.Item(31, 14).Value = -1017141.94   # N31: -1299184.4 -> -1017141.94
$ws.Cells.Item(32, 8).Value = 6500   # H32: 7503.3335 -> 6500
$ws.Cells.Item(32, 9).Value = 10000   # I32: 7503.3335 -> 10000
$ws.Cells.Item(32, 10).Value = 3000   # J32: 0 -> 3000
$ws.Cells.Item(32, 11).Value = 10000   # K32: 7503.3335 -> 10000
$ws.Cells.Item(32, 12).Value = 3000   # L32: 0 -> 3000
$ws.Cells.Item(32, 13).Value = -9684   # M32: -7187.3335 -> -9684
$ws.Cells.Item(32, 14).Value = -3632   # N32: None -> -3632
$ws.Cells.Item(34, 8).Value = 561191.5600000001   # H34: 713937.9399999999 -> 561191.5600000001
$ws.Cells.Item(34, 9).Value = 9965.842000000001   # I34: 12350.134 -> 9965.842000000001
$ws.Cells.Item(34, 10).Value = 1016551.94   # J34: 1298594.4 -> 1016551.94
$ws.Cells.Item(34, 11).Value = 9965.842000000001   # K34: 12350.134 -> 9965.842000000001
$ws.Cells.Item(34, 12).Value = 1016551.94   # L34: 1298594.4 -> 1016551.94
$ws.Cells.Item(34, 13).Value = -9763.842000000001   # M34: -12148.134 -> -9763.842000000001
$ws.Cells.Item(34, 14).Value = -1016955.94   # N34: -1298998.4 -> -1016955.94
$ws.Cells.Item(35, 8).Value = 20000000   # H35: 10000500 -> 20000000
$ws.Cells.Item(35, 9).Value = 0   # I35: 1000 -> 0
$ws.Cells.Item(35, 11).Value = 0   # K35: 1000 -> 0
$ws.Cells.Item(35, 13).Value = ""   # M35: -706 -> (cell removed)
$ws.Cells.Item(56, 8).Value = 19000   # H56: 0 -> 19000
$ws.Cells.Item(56, 10).Value = 19000   # J56: 0 -> 19000
$ws.Cells.Item(56, 12).Value = 19000   # L56: 0 -> 19000
$ws.Cells.Item(56, 14).Value = -20690   # N56: None -> -20690
$ws.Cells.Item(62, 8).Value = 3530.5557   # H62: 3648.125 -> 3530.5557
$ws.Cells.Item(62, 10).Value = 3635.6   # J62: 3897 -> 3635.6
$ws.Cells.Item(62, 12).Value = 3635.6   # L62: 3897 -> 3635.6
$ws.Cells.Item(62, 14).Value = -4883.6   # N62: -5145 -> -4883.6
$ws.Cells.Item(65, 8).Value = 3530.5557   # H65: 3648.125 -> 3530.5557
$ws.Cells.Item(65, 10).Value = 3635.6   # J65: 3897 -> 3635.6
$ws.Cells.Item(65, 12).Value = 18178   # L65: 19485 -> 18178
$ws.Cells.Item(65, 14).Value = -24418   # N65: -25725 -> -24418
$ws.Cells.Item(105, 8).Value = 1179.8462   # H105: 1223.25 -> 1179.8462
$ws.Cells.Item(105, 9).Value = 964   # I105: 997.8889 -> 964
$ws.Cells.Item(105, 11).Value = 964   # K105: 997.8889 -> 964
$ws.Cells.Item(105, 13).Value = 783   # M105: 749.1111 -> 783
$ws.Cells.Item(132, 8).Value = 2407.8064   # H132: 2446.9333 -> 2407.8064
$ws.Cells.Item(132, 9).Value = 2213.6538   # I132: 2248.6667 -> 2213.6538
$ws.Cells.Item(132, 10).Value = 3417.4   # J132: 4231.3335 -> 3417.4
$ws.Cells.Item(132, 11).Value = 6640.9614   # K132: 6746.000100000001 -> 6640.9614
$ws.Cells.Item(132, 12).Value = 10252.2   # L132: 12694.0005 -> 10252.2
$ws.Cells.Item(132, 13).Value = -4110.9614   # M132: -4216.000100000001 -> -4110.9614
$ws.Cells.Item(132, 14).Value = -15312.2   # N132: -17754.0005 -> -15312.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 700   # H14: 975 -> 700
$ws.Cells.Item(14, 9).Value = 700   # I14: 975 -> 700
$ws.Cells.Item(14, 11).Value = 2100   # K14: 2925 -> 2100
$ws.Cells.Item(14, 13).Value = -1927   # M14: -2752 -> -1927
$ws.Cells.Item(37, 8).Value = 79999.5   # H37: 80000 -> 79999.5
$ws.Cells.Item(37, 10).Value = 79999.5   # J37: 80000 -> 79999.5
$ws.Cells.Item(37, 12).Value = 239998.5   # L37: 240000 -> 239998.5
$ws.Cells.Item(37, 14).Value = -240222.5   # N37: -240224 -> -240222.5
$ws.Cells.Item(68, 8).Value = 1629.6171   # H68: 1635.1702 -> 1629.6171
$ws.Cells.Item(68, 10).Value = 1583.8649   # J68: 1590.919 -> 1583.8649
$ws.Cells.Item(68, 12).Value = 4751.5947   # L68: 4772.757000000001 -> 4751.5947
$ws.Cells.Item(68, 14).Value = -6373.5947   # N68: -6394.757000000001 -> -6373.5947
$ws.Cells.Item(71, 8).Value = 1629.6171   # H71: 1635.1702 -> 1629.6171
$ws.Cells.Item(71, 10).Value = 1583.8649   # J71: 1590.919 -> 1583.8649
$ws.Cells.Item(71, 12).Value = 14254.7841   # L71: 14318.271 -> 14254.7841
$ws.Cells.Item(71, 14).Value = -22366.7841   # N71: -22430.271 -> -22366.7841
$ws.Cells.Item(92, 8).Value = 1430784.1   # H92: 436152.44 -> 1430784.1
$ws.Cells.Item(92, 9).Value = 3336966.2   # I92: 1668983.9 -> 3336966.2
$ws.Cells.Item(92, 10).Value = 1147.5   # J92: 1035.4706 -> 1147.5
$ws.Cells.Item(92, 11).Value = 10010898.6   # K92: 5006951.699999999 -> 10010898.6
$ws.Cells.Item(92, 12).Value = 3442.5   # L92: 3106.4118 -> 3442.5
$ws.Cells.Item(92, 13).Value = -10009650.6   # M92: -5005703.699999999 -> -10009650.6
$ws.Cells.Item(92, 14).Value = -5938.5   # N92: -5602.4118 -> -5938.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5161.875   # H70: 5584 -> 5161.875
$ws.Cells.Item(70, 9).Value = 5199.4287   # I70: 5755.25 -> 5199.4287
$ws.Cells.Item(70, 11).Value = 5199.4287   # K70: 5755.25 -> 5199.4287
$ws.Cells.Item(70, 13).Value = -4929.4287   # M70: -5485.25 -> -4929.4287
$ws.Cells.Item(73, 8).Value = 5161.875   # H73: 5584 -> 5161.875
$ws.Cells.Item(73, 9).Value = 5199.4287   # I73: 5755.25 -> 5199.4287
$ws.Cells.Item(73, 11).Value = 5199.4287   # K73: 5755.25 -> 5199.4287
$ws.Cells.Item(73, 13).Value = -4263.4287   # M73: -4819.25 -> -4263.4287
$ws.Cells.Item(80, 8).Value = 23592.072   # H80: 24184.54 -> 23592.072
$ws.Cells.Item(80, 9).Value = 20348.75   # I80: 20985.715 -> 20348.75
$ws.Cells.Item(80, 11).Value = 20348.75   # K80: 20985.715 -> 20348.75
$ws.Cells.Item(80, 13).Value = -19350.75   # M80: -19987.715 -> -19350.75
$ws.Cells.Item(83, 8).Value = 23592.072   # H83: 24184.54 -> 23592.072
$ws.Cells.Item(83, 9).Value = 20348.75   # I83: 20985.715 -> 20348.75
$ws.Cells.Item(83, 11).Value = 101743.75   # K83: 104928.575 -> 101743.75
$ws.Cells.Item(83, 13).Value = -96751.75   # M83: -99936.575 -> -96751.75
$ws.Cells.Item(97, 8).Value = 1548.6786   # H97: 1579.8148 -> 1548.6786
$ws.Cells.Item(97, 9).Value = 1570.9231   # I97: 1605.44 -> 1570.9231
$ws.Cells.Item(97, 11).Value = 1570.9231   # K97: 1605.44 -> 1570.9231
$ws.Cells.Item(97, 13).Value = -1074.9231   # M97: -1109.44 -> -1074.9231

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3039.2307   # H16: 3088.84 -> 3039.2307
$ws.Cells.Item(16, 9).Value = 3208.0557   # I16: 3290.9412 -> 3208.0557
$ws.Cells.Item(16, 11).Value = 3208.0557   # K16: 3290.9412 -> 3208.0557
$ws.Cells.Item(16, 13).Value = -3038.0557   # M16: -3120.9412 -> -3038.0557
$ws.Cells.Item(22, 8).Value = 6714.8887   # H22: 12777 -> 6714.8887
$ws.Cells.Item(22, 9).Value = 6859.875   # I22: 19999 -> 6859.875
$ws.Cells.Item(22, 11).Value = 6859.875   # K22: 19999 -> 6859.875
$ws.Cells.Item(22, 13).Value = -6564.875   # M22: -19704 -> -6564.875
$ws.Cells.Item(27, 8).Value = 6714.8887   # H27: 12777 -> 6714.8887
$ws.Cells.Item(27, 9).Value = 6859.875   # I27: 19999 -> 6859.875
$ws.Cells.Item(27, 11).Value = 6859.875   # K27: 19999 -> 6859.875
$ws.Cells.Item(27, 13).Value = -6752.875   # M27: -19892 -> -6752.875
$ws.Cells.Item(40, 8).Value = 4494.8887   # H40: 4833.3335 -> 4494.8887
$ws.Cells.Item(40, 9).Value = 3863.5   # I40: 4000 -> 3863.5
$ws.Cells.Item(40, 11).Value = 3863.5   # K40: 4000 -> 3863.5
$ws.Cells.Item(40, 13).Value = -3727.5   # M40: -3864 -> -3727.5
$ws.Cells.Item(55, 8).Value = 62500270   # H55: 50000228 -> 62500270
$ws.Cells.Item(55, 9).Value = 111111360   # I55: 76923260 -> 111111360
$ws.Cells.Item(55, 11).Value = 111111360   # K55: 76923260 -> 111111360
$ws.Cells.Item(55, 13).Value = -111111187   # M55: -76923087 -> -111111187
$ws.Cells.Item(122, 8).Value = 4446.972   # H122: 4725.0605 -> 4446.972
$ws.Cells.Item(122, 9).Value = 4158.9375   # I122: 4445.5864 -> 4158.9375
$ws.Cells.Item(122, 11).Value = 12476.8125   # K122: 13336.7592 -> 12476.8125
$ws.Cells.Item(122, 13).Value = -10026.8125   # M122: -10886.7592 -> -10026.8125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2249   # H96: 2274 -> 2249
$ws.Cells.Item(96, 9).Value = 2235.5789   # I96: 2267.158 -> 2235.5789
$ws.Cells.Item(96, 11).Value = 2235.5789   # K96: 2267.158 -> 2235.5789
$ws.Cells.Item(96, 13).Value = -862.5789   # M96: -894.1579999999999 -> -862.5789
$ws.Cells.Item(132, 8).Value = 427007.22   # H132: 418119.62 -> 427007.22
$ws.Cells.Item(132, 9).Value = 1437.7073   # I132: 1413.0714 -> 1437.7073
$ws.Cells.Item(132, 11).Value = 4313.1219   # K132: 4239.2142 -> 4313.1219
$ws.Cells.Item(132, 13).Value = -1783.1219   # M132: -1709.2142 -> -1783.1219
